# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) sheet (4th worksheet) gets three new trailing columns —
# date, legislator_name, legislator_id — mirroring the same three columns
# already appended to the other property sheets for this legislator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$xlPasteFormats = -4122

# --- Header row (row 1): new column headers, styled like the existing ones ---
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows (2-4): copy each row's existing formatting onto its new cells ---
$ws.Range("G2").Copy()
$ws.Range("H2:J2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("G3").Copy()
$ws.Range("H3:J3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("G4").Copy()
$ws.Range("H4:J4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Force the date column to text first so "2013-12-31" is stored as a literal
# string rather than auto-converted into a date serial number.
$ws.Range("H2:H4").NumberFormat = "@"

$ws.Range("H2").Value = "2013-12-31"
$ws.Range("I2").Value = "楊麗環"
$ws.Range("J2").Value = 960

$ws.Range("H3").Value = "2013-12-31"
$ws.Range("I3").Value = "楊麗環"
$ws.Range("J3").Value = 960

$ws.Range("H4").Value = "2013-12-31"
$ws.Range("I4").Value = "楊麗環"
$ws.Range("J4").Value = 960

# Re-apply the original (General-format) look to the date cells now that the
# text value is safely stored, so they match the plain styling of the rest
# of the row instead of keeping a distinct "text" number format.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
